$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 7) mirroring the existing rows' structure.
$row = 7

$ws.Cells.Item($row, 1).Value = 42604.891469907408
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = "Noun"

$ws.Cells.Item($row, 3).Value = 2011
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 2
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
